# Bugfixed evaluation and simulated rt_data for components
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..53 (date serial in column A, value in column B)
$data = @(
@(39400,-3.518004331294975),
@(39583,4.488226581078521),
@(39765,-3.397482474306742),
@(39948,-2.589909567620779),
@(40130,1.250154354816019),
@(40310,-0.9550786931962989),
@(40494,6.853351083518234),
@(40676,-2.590749655123886),
@(40862,-0.66979119045169),
@(41044,-1.346695203931787),
@(41228,1.527423145498517),
@(41409,-2.087743866862297),
@(41592,1.910794643105291),
@(41774,3.640709258360772),
@(41957,-3.894663176624903),
@(42137,1.325647082233132),
@(42321,-1.301240742391741),
@(42503,1.969890649492385),
@(42689,0.279053692590864),
@(42867,0.8237845450787233),
@(43053,0.494904046794133),
@(43145,-0.2521030831923952),
@(43235,2.072093850483441),
@(43326,1.598757058788109),
@(43418,0.9),
@(43510,0.7186276419777187),
@(43600,1.9),
@(43691,2.478732535894252),
@(43783,-1.062312873792649),
@(43875,0.3597775049439349),
@(43966,4.1),
@(44068,-4.222133416425748),
@(44159,-1.967214790670653),
@(44251,1.836112747457122),
@(44341,1.142507241736567),
@(44432,0.2778921468196529),
@(44525,-2.336211685220064),
@(44617,0.03642954902896633),
@(44706,4.623062809051376),
@(44798,-3.361711150938063),
@(44890,-1.432989117992093),
@(44981,-2.851410521811431),
@(45071,3.933557394002079),
@(45163,0.2329359358994338),
@(45254,0.4324506349439758),
@(45345,-1.721673722292465),
@(45436,2.712662502476618),
@(45534,-2.00812900028761),
@(45618,-0.346469331090276),
@(45713,1.043712238889796),
@(45800,0.5275848603473889),
@(45891,-2.073458279789094)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Remove the now-obsolete trailing rows (previously up to row 73)
$lastOldRow = 73
$firstRowToDelete = $startRow + $data.Count
if ($lastOldRow -ge $firstRowToDelete) {
    $deleteRange = $ws.Range("A" + $firstRowToDelete + ":A" + $lastOldRow)
    $deleteRange.EntireRow.Delete()
}
